$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.020335
$ws.Range("H2").Value = 0.061005
$ws.Range("I2").Value = 0.009804808687698561
$ws.Range("J2").Value = 0.009804808687698559
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3407069999999999
$ws.Range("N2").Value = 1.022121
$ws.Range("O2").Value = 0.1055965976712818
$ws.Range("P2").Value = 0.1055965976712818
$ws.Range("Q2").Value = 0.006928276844999998
$ws.Range("R2").Value = 0.06235449160499999
$ws.Range("S2").Value = 0.001035354438238793
$ws.Range("T2").Value = 0.001035354438238793
$ws.Range("G3").Value = 0.020335
$ws.Range("H3").Value = 0.061005
$ws.Range("I3").Value = 0.009804808687698561
$ws.Range("J3").Value = 0.009804808687698559
$ws.Range("O3").Value = 0.3782153560188308
$ws.Range("P3").Value = 0.3782153560188308
$ws.Range("Q3").Value = 0.02481501062833333
$ws.Range("R3").Value = 0.223335095655
$ws.Range("S3").Value = 0.003708329208514436
$ws.Range("T3").Value = 0.003708329208514436
$ws.Range("G4").Value = 0.020335
$ws.Range("H4").Value = 0.061005
$ws.Range("I4").Value = 0.009804808687698561
$ws.Range("J4").Value = 0.009804808687698559
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5161880463098875
$ws.Range("P4").Value = 0.5161880463098875
$ws.Range("Q4").Value = 0.03386750868666666
$ws.Range("R4").Value = 0.3048075781799999
$ws.Range("S4").Value = 0.005061125040945332
$ws.Range("T4").Value = 0.005061125040945331
$ws.Range("I5").Value = 0.1486140913768632
$ws.Range("J5").Value = 0.1486140913768632
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3407069999999999
$ws.Range("N5").Value = 1.022121
$ws.Range("O5").Value = 0.1055965976712818
$ws.Range("P5").Value = 0.1055965976712818
$ws.Range("Q5").Value = 0.105013733661
$ws.Range("R5").Value = 0.9451236029489999
$ws.Range("S5").Value = 0.01569314241540574
$ws.Range("T5").Value = 0.01569314241540573
$ws.Range("I6").Value = 0.1486140913768632
$ws.Range("J6").Value = 0.1486140913768632
$ws.Range("O6").Value = 0.3782153560188308
$ws.Range("P6").Value = 0.3782153560188308
$ws.Range("S6").Value = 0.05620813147951538
$ws.Range("T6").Value = 0.05620813147951537
$ws.Range("I7").Value = 0.1486140913768632
$ws.Range("J7").Value = 0.1486140913768632
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5161880463098875
$ws.Range("P7").Value = 0.5161880463098875
$ws.Range("Q7").Value = 0.513338831076
$ws.Range("R7").Value = 4.620049479684
$ws.Range("S7").Value = 0.07671281748194213
$ws.Range("T7").Value = 0.07671281748194211
$ws.Range("H8").Value = 5.236273000000001
$ws.Range("I8").Value = 0.8415810999354383
$ws.Range("J8").Value = 0.8415810999354382
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3407069999999999
$ws.Range("N8").Value = 1.022121
$ws.Range("O8").Value = 0.1055965976712818
$ws.Range("P8").Value = 0.1055965976712818
$ws.Range("Q8").Value = 0.5946782883369999
$ws.Range("R8").Value = 5.352104595033
$ws.Range("S8").Value = 0.08886810081763727
$ws.Range("T8").Value = 0.08886810081763725
$ws.Range("H9").Value = 5.236273000000001
$ws.Range("I9").Value = 0.8415810999354383
$ws.Range("J9").Value = 0.8415810999354382
$ws.Range("O9").Value = 0.3782153560188308
$ws.Range("P9").Value = 0.3782153560188308
$ws.Range("S9").Value = 0.318298895330801
$ws.Range("T9").Value = 0.318298895330801
$ws.Range("H10").Value = 5.236273000000001
$ws.Range("I10").Value = 0.8415810999354383
$ws.Range("J10").Value = 0.8415810999354382
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5161880463098875
$ws.Range("P10").Value = 0.5161880463098875
$ws.Range("Q10").Value = 2.906966991447555
$ws.Range("S10").Value = 0.434414103787
$ws.Range("T10").Value = 0.434414103787
